# Revised data files to make them R-friendly
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header from "Sale" to "sale_amount"
$ws.Range("A1").Value = "sale_amount"

# Strip the currency formatting from the data column, switching it to a
# plain two-decimal numeric format (easier to read into R cleanly).
$ws.Range("A1:A32").NumberFormat = "0.00"

# Select the whole column, as was left selected in the saved file.
$ws.Columns.Item(1).Select()
